$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 3 and row 4: A, B, D, E, F, G, H
$cols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $val3 = $cell3.Value2
    $val4 = $cell4.Value2
    $cell3.Value = $val4
    $cell4.Value = $val3
}
